# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# on all three sheets, then resize the affected "Status" columns the way a
# report-refresh pass would (shrinking them to fit the shorter text).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
if ($overview.Range("E2").Text -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Text -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($overview.Range("E3").Text -eq $oldStatus) { $overview.Range("E3").Value = $newStatus }
if ($overview.Range("F3").Text -eq $oldStatus) { $overview.Range("F3").Value = $newStatus }

# --- zh-cn / de-de detail sheets: Status column (C) ---
if ($zhcn.Range("C2").Text -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($zhcn.Range("C3").Text -eq $oldStatus) { $zhcn.Range("C3").Value = $newStatus }

if ($dede.Range("C2").Text -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }
if ($dede.Range("C3").Text -eq $oldStatus) { $dede.Range("C3").Value = $newStatus }

# --- Re-fit the now-narrower status columns so the sheet matches the
#     refreshed report layout. ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
